$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4287.976
$ws.Range("I15").Value = 4287.976
$ws.Range("K15").Value = 12863.928
$ws.Range("M15").Value = -12694.928
$ws.Range("H17").Value = 605.375
$ws.Range("J17").Value = 605.375
$ws.Range("L17").Value = 1816.125
$ws.Range("N17").Value = -2152.125
$ws.Range("H39").Value = 586.2222
$ws.Range("I39").Value = 122.14286
$ws.Range("J39").Value = 881.5454999999999
$ws.Range("K39").Value = 366.42858
$ws.Range("L39").Value = 2644.6365
$ws.Range("M39").Value = -70.42858000000001
$ws.Range("N39").Value = -3236.6365
$ws.Range("H135").Value = 60074.176
$ws.Range("I135").Value = 60074.176
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 540667.584
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -538132.584
$ws.Range("N135").ClearContents()
$ws.Range("H141").Value = 1018.2727
$ws.Range("I141").Value = 620.1
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 1860.3
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = 3319.7
$ws.Range("N141").Value = -25360

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2105.6428
$ws.Range("I2").Value = 2117.9
$ws.Range("J2").Value = 2075
$ws.Range("K2").Value = 2117.9
$ws.Range("L2").Value = 2075
$ws.Range("M2").Value = -2004.9
$ws.Range("N2").Value = -2301
$ws.Range("H32").Value = 16505.023
$ws.Range("I32").Value = 19263.857
$ws.Range("J32").Value = 8228.522999999999
$ws.Range("K32").Value = 19263.857
$ws.Range("L32").Value = 8228.522999999999
$ws.Range("M32").Value = -18976.857
$ws.Range("N32").Value = -8802.522999999999
$ws.Range("H61").Value = 34552564
$ws.Range("I61").Value = 45500628
$ws.Range("J61").Value = 144360.72
$ws.Range("K61").Value = 45500628
$ws.Range("L61").Value = 144360.72
$ws.Range("M61").Value = -45500416
$ws.Range("N61").Value = -144784.72
$ws.Range("H74").Value = 3707382.8
$ws.Range("I74").Value = 4827870
$ws.Range("J74").Value = 65798.375
$ws.Range("K74").Value = 4827870
$ws.Range("L74").Value = 65798.375
$ws.Range("M74").Value = -4826996
$ws.Range("N74").Value = -67546.375
$ws.Range("H77").Value = 3707382.8
$ws.Range("I77").Value = 4827870
$ws.Range("J77").Value = 65798.375
$ws.Range("K77").Value = 24139350
$ws.Range("L77").Value = 328991.875
$ws.Range("M77").Value = -24134982
$ws.Range("N77").Value = -337727.875
$ws.Range("H116").Value = 2105.6428
$ws.Range("I116").Value = 2117.9
$ws.Range("J116").Value = 2075
$ws.Range("K116").Value = 2117.9
$ws.Range("L116").Value = 2075
$ws.Range("M116").Value = 176.0999999999999
$ws.Range("N116").Value = -6663
$ws.Range("H122").Value = 2850688.8
$ws.Range("I122").Value = 1715.3334
$ws.Range("J122").Value = 37038372
$ws.Range("K122").Value = 5146.0002
$ws.Range("L122").Value = 111115116
$ws.Range("M122").Value = -2696.0002
$ws.Range("N122").Value = -111120016
$ws.Range("H132").Value = 52415.562
$ws.Range("I132").Value = 35380.656
$ws.Range("J132").Value = 93583.25
$ws.Range("K132").Value = 106141.968
$ws.Range("L132").Value = 280749.75
$ws.Range("M132").Value = -103611.968
$ws.Range("N132").Value = -285809.75
$ws.Range("H133").Value = 49800
$ws.Range("J133").Value = 49800
$ws.Range("L133").Value = 49800
$ws.Range("N133").Value = -54860
$ws.Range("H136").Value = 34552564
$ws.Range("I136").Value = 45500628
$ws.Range("J136").Value = 144360.72
$ws.Range("K136").Value = 136501884
$ws.Range("L136").Value = 433082.16
$ws.Range("M136").Value = -136499334
$ws.Range("N136").Value = -438182.16

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2105.6428
$ws.Range("I3").Value = 2117.9
$ws.Range("J3").Value = 2075
$ws.Range("K3").Value = 2117.9
$ws.Range("L3").Value = 2075
$ws.Range("M3").Value = -2003.9
$ws.Range("N3").Value = -2303
$ws.Range("H30").Value = 12000
$ws.Range("J30").Value = 12000
$ws.Range("L30").Value = 12000
$ws.Range("N30").Value = -12250
$ws.Range("H99").Value = 1121.3043
$ws.Range("I99").Value = 1105.2941
$ws.Range("K99").Value = 1105.2941
$ws.Range("M99").Value = 392.7058999999999
$ws.Range("H107").Value = 2403.3076
$ws.Range("I107").Value = 2737
$ws.Range("J107").Value = 2117.2856
$ws.Range("K107").Value = 2737
$ws.Range("L107").Value = 2117.2856
$ws.Range("M107").Value = -817
$ws.Range("N107").Value = -5957.2856
$ws.Range("H123").Value = 23583.166
$ws.Range("J123").Value = 23583.166
$ws.Range("L123").Value = 23583.166
$ws.Range("N123").Value = -33383.166

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 1000
$ws.Range("K2").Value = 1000
$ws.Range("M2").Value = -887
$ws.Range("H58").Value = 31252342
$ws.Range("I58").Value = 43480344
$ws.Range("J58").Value = 3000.111
$ws.Range("K58").Value = 43480344
$ws.Range("L58").Value = 3000.111
$ws.Range("M58").Value = -43480141
$ws.Range("N58").Value = -3406.111
$ws.Range("H94").Value = 5300.273
$ws.Range("I94").Value = 12130.75
$ws.Range("J94").Value = 1397.1428
$ws.Range("K94").Value = 12130.75
$ws.Range("L94").Value = 1397.1428
$ws.Range("M94").Value = -11679.75
$ws.Range("N94").Value = -2299.1428
$ws.Range("H99").Value = 4274
$ws.Range("I99").Value = 4165.1113
$ws.Range("J99").Value = 4372
$ws.Range("K99").Value = 4165.1113
$ws.Range("L99").Value = 4372
$ws.Range("M99").Value = -2667.1113
$ws.Range("N99").Value = -7368
$ws.Range("H122").Value = 2242.7058
$ws.Range("I122").Value = 1655.6364
$ws.Range("J122").Value = 3319
$ws.Range("K122").Value = 4966.9092
$ws.Range("L122").Value = 9957
$ws.Range("M122").Value = -2516.9092
$ws.Range("N122").Value = -14857
$ws.Range("H126").Value = 4274
$ws.Range("I126").Value = 4165.1113
$ws.Range("J126").Value = 4372
$ws.Range("K126").Value = 12495.3339
$ws.Range("L126").Value = 13116
$ws.Range("M126").Value = -10025.3339
$ws.Range("N126").Value = -18056
$ws.Range("H134").Value = 41025.242
$ws.Range("I134").Value = 2526.8696
$ws.Range("J134").Value = 188602.33
$ws.Range("K134").Value = 7580.6088
$ws.Range("L134").Value = 565806.99
$ws.Range("M134").Value = -5045.6088
$ws.Range("N134").Value = -570876.99
$ws.Range("H136").Value = 31252342
$ws.Range("I136").Value = 43480344
$ws.Range("J136").Value = 3000.111
$ws.Range("K136").Value = 130441032
$ws.Range("L136").Value = 9000.332999999999
$ws.Range("M136").Value = -130438482
$ws.Range("N136").Value = -14100.333

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 2757.1428
$ws.Range("I118").Value = 400
$ws.Range("J118").Value = 3700
$ws.Range("K118").Value = 1200
$ws.Range("L118").Value = 11100
$ws.Range("M118").Value = 43
$ws.Range("N118").Value = -13586
$ws.Range("H131").Value = 1216.4182
$ws.Range("J131").Value = 1274.06
$ws.Range("L131").Value = 3822.18
$ws.Range("N131").Value = -13902.18

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 6000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 6000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 6000
$ws.Range("N23").Value = -6446
$ws.Range("M23").ClearContents()
$ws.Range("H102").Value = 1114.091
$ws.Range("I102").Value = 906
$ws.Range("J102").Value = 1160.3334
$ws.Range("K102").Value = 906
$ws.Range("L102").Value = 1160.3334
$ws.Range("M102").Value = 716
$ws.Range("N102").Value = -4404.3334

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3498
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 3498
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 3498
$ws.Range("N7").Value = -3722
$ws.Range("M7").ClearContents()
$ws.Range("H40").Value = 2500
$ws.Range("I40").Value = 2500
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 2500
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -2364
$ws.Range("N40").ClearContents()
$ws.Range("H68").Value = 1598.2609
$ws.Range("I68").Value = 1603
$ws.Range("J68").Value = 1566.6666
$ws.Range("K68").Value = 1603
$ws.Range("L68").Value = 1566.6666
$ws.Range("M68").Value = -854
$ws.Range("N68").Value = -3064.6666
$ws.Range("H71").Value = 1598.2609
$ws.Range("I71").Value = 1603
$ws.Range("J71").Value = 1566.6666
$ws.Range("K71").Value = 8015
$ws.Range("L71").Value = 7833.333000000001
$ws.Range("M71").Value = -4271
$ws.Range("N71").Value = -15321.333
$ws.Range("H93").Value = 2273.2727
$ws.Range("I93").Value = 2072.2856
$ws.Range("K93").Value = 2072.2856
$ws.Range("M93").Value = -824.2856000000002
$ws.Range("H100").Value = 1799.8572
$ws.Range("I100").Value = 1559.6
$ws.Range("J100").Value = 1933.3334
$ws.Range("K100").Value = 1559.6
$ws.Range("L100").Value = 1933.3334
$ws.Range("M100").Value = -1018.6
$ws.Range("N100").Value = -3015.3334
$ws.Range("H118").Value = 32000
$ws.Range("J118").Value = 32000
$ws.Range("L118").Value = 32000
$ws.Range("N118").Value = -35314
$ws.Range("H126").Value = 3498
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3498
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 10494
$ws.Range("N126").Value = -15434
$ws.Range("M126").ClearContents()
$ws.Range("H136").Value = 41228
$ws.Range("I136").Value = 24376.232
$ws.Range("J136").Value = 131806.25
$ws.Range("K136").Value = 73128.696
$ws.Range("L136").Value = 395418.75
$ws.Range("M136").Value = -70578.696
$ws.Range("N136").Value = -400518.75

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1561
$ws.Range("I126").Value = 1356.1428
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 4068.4284
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -1598.4284
$ws.Range("N126").Value = -10340
$ws.Range("H129").Value = 29620
$ws.Range("J129").Value = 29620
$ws.Range("L129").Value = 29620
$ws.Range("N129").Value = -39620
